$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly data between row 2 (newer) and row 3 (older):
# Date
$ws.Range("D2").Value = 44291
$ws.Range("D3").Value = 44421

# Volumen
$ws.Range("M2").Value = 15
$ws.Range("M3").Value = 30

# Precio minimo
$ws.Range("N2").Value = 23000
$ws.Range("N3").Value = 24000

# Precio maximo
$ws.Range("O2").Value = 23000
$ws.Range("O3").Value = 24000

# Precio promedio ponderado
$ws.Range("P2").Value = 23000
$ws.Range("P3").Value = 24000

# Precio $/Kg
$ws.Range("S2").Value = 1150
$ws.Range("S3").Value = 1200
